$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing EntryScene values for the revive-in-dungeon fix
$ws.Range("E5").Value = 13010007
$ws.Range("E6").Value = 13010004
$ws.Range("E7").Value = 13010005

# Update selection to match author's last cursor position
$ws.Range("E7").Select()

# Reset the workbook's "Light 1" theme color back to plain white
$theme = $wb.Theme
$scheme = $theme.ThemeColorScheme
$scheme.Colors(2).RGB = 0xFFFFFF
